$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1 - mirror the style of the existing header (H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data values for I2:J55
$iValues = @(8,9,9,7,9,4,9,9,8,6,6,6,8,6,7,6,6,9,6,7,9,5,7,6,6,4,3,5,8,5,6,8,6,7,8,6,7,9,7,8,6,8,7,7,7,9,5,4,5,7,7,7,4,2)
$jValues = @(8,9,10,8,9,5,9,9,8,7,7,7,8,6,7,6,6,9,7,8,9,5,7,6,6,5,3,6,8,5,6,8,7,7,8,6,7,9,7,8,6,8,7,7,7,9,5,4,5,7,7,7,4,2)

for ($n = 0; $n -lt $iValues.Length; $n++) {
    $row = $n + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$n]
    $ws.Cells.Item($row, 10).Value = $jValues[$n]
}
